# Updates cryptos worksheet: refresh Price (D) and Volume(1h) (E) columns
# with the latest scraped figures, and fix the InjectiveProtocol/Hedera
# row ordering (rows 46-47 swapped back with refreshed data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.435.51"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "2.328.46"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'542.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.38%  "
$ws.Range("D6").Value = "'134.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "2.363.51"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").Value = "'0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").Value = "'5.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").Value = "2.782.13"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "'23.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "57.607.11"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "2.348.65"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").Value = "'337.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D20").Value = "'10.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").Value = "'4.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").Value = "'6.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.25%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'61.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").Value = "'0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("D26").Value = "'8.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'1.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.51%  "
$ws.Range("E29").Value = "  +5.91%  "
$ws.Range("D30").Value = "'171.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "0.0₃0737"
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("D32").Value = "'6.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").Value = "'18.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("E34").Value = "  +15.76%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'0.990"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").Value = "'4.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.23%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +5.08%  "
$ws.Range("D40").Value = "'39.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'148.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").Value = "'285.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").Value = "'0.0933"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0507"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'19.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.61%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").Value = "'17.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").Value = "'0.380"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
